$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:D8").Select()
Write-Output "ActiveCell before:"
Write-Output $excel.ActiveCell.Address()
$excel.ActiveWindow.RangeSelection.Item(2,1).Activate()
Write-Output "ActiveCell after:"
Write-Output $excel.ActiveCell.Address()
